$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(39, 41, 42, 43, 44, 55, 56, 58)
foreach ($r in $rows) {
    $ws.Range("H$r").Value = -1
}
